$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the state-function matrix: several EFT flags flip from 1 to 0
# for the listed Func rows (B..O are EFT_INFO..EFT_FINISH columns).
$ws.Range("E2:N2").Value = 0
$ws.Range("C3:D3").Value = 0
$ws.Range("F3:N3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4:N4").Value = 0
$ws.Range("C5:D5").Value = 0
$ws.Range("F5:N5").Value = 0
$ws.Range("C6:N6").Value = 0
$ws.Range("C7:N7").Value = 0
$ws.Range("C8:N8").Value = 0
$ws.Range("C9:N9").Value = 0
$ws.Range("C10:N10").Value = 0
$ws.Range("C11:N11").Value = 0
$ws.Range("C12:N12").Value = 0
$ws.Range("C13:N13").Value = 0

# Move the active selection to match the saved workbook state.
$ws.Range("F11").Select()

$wb.Save()
